$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ER)
$ws.Range("B2").Value = 71.66704759788001
$ws.Range("C2").Value = 95.96322044903999
$ws.Range("D2").Value = 84.28312065108
$ws.Range("E2").Value = 5605.17276423382
$ws.Range("F2").Value = 7739.434472504409

# Row 3 (SF)
$ws.Range("B3").Value = 43.63070991305
$ws.Range("C3").Value = 110.11073297724
$ws.Range("D3").Value = 99.17654561585
$ws.Range("E3").Value = 4772.494032335781
$ws.Range("F3").Value = 24908.75805290864

# Row 4 (WA)
$ws.Range("B4").Value = 47.09307318523
$ws.Range("C4").Value = 78.00284813434999
$ws.Range("D4").Value = 47.94975442803
$ws.Range("E4").Value = 4365.34005233962
$ws.Range("F4").Value = 5246.26399374735
